$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so leftover legacy styling (old gray fill, etc.)
# does not linger on cells we are about to rewrite.
$ws.Cells.Clear()

# --- Header row ---
$headers = @("Номер пользователя", "Логин (Аккаунт пользователя на TimeWeb)", "Стоимость тарифа", "Последняя дата оплаты")
for ($col = 1; $col -le 4; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Data rows: phone number, login, tariff cost, last payment date (serial) ---
$data = @(
    @("8(908)-642-00-01", "co75007", 270,  45361),
    @("8(908)-642-00-02", "co75008", 480,  45362),
    @("8(908)-642-00-03", "co75009", 690,  45363),
    @("8(908)-642-00-04", "co75010", 900,  45364),
    @("8(908)-642-00-05", "co75011", 1110, 45365),
    @("8(908)-642-00-06", "co75012", 1320, 45366),
    @("8(908)-642-00-07", "co75013", 1530, 45367),
    @("8(908)-642-00-08", "co75014", 1740, 45368),
    @("8(908)-642-00-09", "co75015", 1950, 45369),
    @("8(908)-642-00-10", "co75016", 2160, 45370)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $record = $data[$i]
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
}

# --- Formatting: font size 10, centered, for the whole table ---
$full = $ws.Range("A1:D11")
$full.Font.Size = 10
$full.HorizontalAlignment = -4108

# --- Date column formatting ---
$dateRange = $ws.Range("D2:D11")
$dateRange.NumberFormat = "dd/mm/yyyy"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 27.08
$ws.Columns.Item(2).ColumnWidth = 44.65
$ws.Columns.Item(3).ColumnWidth = 20.36
$ws.Columns.Item(4).ColumnWidth = 21.51
